$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 1: title text box ("Rectangle 4") - merge the "Lecture " / "six: "
# runs into a single "Lecture 6: " run (the rest of the sentence is left
# untouched).
# ---------------------------------------------------------------------------
$titleSlide = $p.Slides.Item(1)
$titleShape = $titleSlide.Shapes.Item(5)
$titleRange = $titleShape.TextFrame.TextRange
$titleLead = $titleRange.Characters(1, 13)
$titleLead.Text = "Lecture 6: "

# ---------------------------------------------------------------------------
# Slide 38 (last slide): "TextBox 2" placeholder that used to contain the
# "<A synopsis of first CUDA lecture>" / "Need to see lecture before I can
# fill this in!" placeholder paragraphs. Replace all of that with the real
# synopsis sentence and collapse the now-empty trailing paragraphs.
# ---------------------------------------------------------------------------
$lastSlide = $p.Slides.Item(38)
$synopsisShape = $lastSlide.Shapes.Item(2)
$synopsisRange = $synopsisShape.TextFrame.TextRange

# Replace everything up to (but not including) the very last paragraph
# mark with the new sentence - this keeps the formatting of the final
# (empty) paragraph mark intact.
$newSentence = "We shall look at CUDA " + [char]0x2013 + " a way to program GPUs"
$body = $synopsisRange.Characters(1, $synopsisRange.Length - 1)
$body.Text = $newSentence

# Drop the now-redundant empty paragraphs that used to separate the two
# placeholder sentences, leaving a single trailing (empty) paragraph mark
# merged right after the new sentence.
$refreshed = $synopsisShape.TextFrame.TextRange
$gapStart = $newSentence.Length + 1
$gapLen = ($refreshed.Length - $newSentence.Length) + 1
$trailing = $refreshed.Characters($gapStart, $gapLen)
$trailing.Delete()
